# Test Data Added for Slovakia market
# Duplicate the "Portugal" sheet (last sheet in the workbook) to create a new
# sheet for the Slovakia market, positioned right after Portugal, then update
# the market name / ticket-reference cells on the new sheet.

$wb = $excel.ActiveWorkbook

$portugal = $wb.Worksheets.Item("Portugal")

# Copy Portugal and place the copy immediately after it.
$portugal.Copy([System.Reflection.Missing]::Value, $portugal) | Out-Null

$slovakia = $wb.Worksheets.Item($portugal.Index + 1)
$slovakia.Name = "Slovakia"

# Update the market name and the JIRA/ticket reference on the new sheet.
$slovakia.Range("B2").Value = "Slovakia Market"
$slovakia.Range("B4").Value = "NGC-2930/T3222"

# The copied rows 3:5 inherited an explicit (wrapped-text) row height from
# Portugal; auto-fit them back down to the sheet's default height.
$slovakia.Rows("3:5").AutoFit() | Out-Null

# Portugal is no longer the active sheet, so select the whole sheet on it
# before moving away (matches its post-edit "not focused" selection state).
$portugal.Activate() | Out-Null
$portugal.Cells.Select() | Out-Null

# Leave the new sheet active, with B4 selected (mirrors where the last edit
# was made) and make it the visible/active tab.
$slovakia.Activate() | Out-Null
$slovakia.Range("B4").Select() | Out-Null
